$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 0.001
$ws.Range("K10").Value = 477
$ws.Range("L10").Value = 0.00159
